$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape id=2 "Rectangle 8": <<interface>>AddressBookStorage -> <<interface>>ThaneParkStorage
$sh2 = $s.Shapes.Item(2)
$sh2.TextFrame.TextRange.Characters(15, 18).Text = "ThaneParkStorage"

# Shape id=50 "Rectangle 8": XmlAddressBook / Storage -> XmlThanePark / Storage
$sh50 = $s.Shapes.Item(13)
$sh50.TextFrame.TextRange.Characters(1, 14).Text = "XmlThanePark"

# Shape id=66 "Rectangle 8": XmlSerializable / AddressBook -> XmlSerializable / ThanePark
$sh66 = $s.Shapes.Item(20)
$sh66.TextFrame.TextRange.Characters(17, 11).Text = "ThanePark"

# Shape id=74 "Rectangle 8": XmlAdaptedPerson -> XmlAdaptedRide
$sh74 = $s.Shapes.Item(23)
$sh74.TextFrame.TextRange.Characters(1, 16).Text = "XmlAdaptedRide"
